$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Description"

$ws.Range("A2").Value = "Move Robot2 to location (2, 8) and remove the toolkit."
$ws.Range("A3").Value = "Move Robot26 to location (11, 4) and remove the liquid spill."
$ws.Range("A4").Value = "Move Robot42 to location (9, 5) and remove the large debris."
$ws.Range("A5").Value = "Move Robot48 to location (5, 6) and remove the dust."
$ws.Range("A6").Value = "Move Robot31 to location (9, 4) and remove the grass."
$ws.Range("A7").Value = "Move Robot8 to location (8, 12) and remove the small debris."
$ws.Range("A8").Value = "Move Robot23 to location (11, 1) and remove the vehicle."
$ws.Range("A9").Value = "Move Robot23 to location (12, 10) and remove the construction materials."
$ws.Range("A10").Value = "Move Robot14 to location (7, 11) and remove the tree branches."
$ws.Range("A11").Value = "Move Robot15 to location (5, 3) and remove the screws."
